$d = $word.ActiveDocument

$d.Content.Find.Execute("1399.33", $true, $false, $false, $false, $false,
                         $true, 1, $false, "151.86", 2)
$d.Content.Find.Execute("1119.46", $true, $false, $false, $false, $false,
                         $true, 1, $false, "121.49", 2)
$d.Content.Find.Execute("1679.20", $true, $false, $false, $false, $false,
                         $true, 1, $false, "182.23", 2)
$d.Content.Find.Execute("GAMMA(100.00, 13.99)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "GAMMA(100.00, 1.52)", 2)
